# Add the "data_dictionary" worksheet after the existing "metro_budget" sheet.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "data_dictionary"

# Column A / B headers + definitions (data dictionary content).
$labels = @(
  "Department",
  "FY17_Budget",
  "FY17_Actual",
  "FY17_diff",
  "FY18_Budget",
  "FY18_Actual",
  "FY18_diff",
  "FY19_Budget",
  "FY19_Actual",
  "FY19_diff"
)

$descriptions = @(
  "Metro Nashville government department name",
  "Department budget for fiscal year 2017",
  "Actual spending for fiscal year 2017",
  "Actual spending amount - budget amount for fiscal year 2017",
  "Department budget for fiscal year 2018",
  "Actual spending for fiscal year 2018",
  "Actual spending amount - budget amount for fiscal year 2018",
  "Department budget for fiscal year 2019",
  "Actual spending for fiscal year 2019",
  "Actual spending amount - budget amount for fiscal year 2019"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $labels[$i]
    $ws2.Cells.Item($row, 2).Value = $descriptions[$i]
}

# Bold black font for column A.
$colA = $ws2.Range("A1:A10")
$colA.Font.Bold = $true
$colA.Font.Color = 0

# Arial 10pt for column B.
$colB = $ws2.Range("B1:B10")
$colB.Font.Name = "Arial"
$colB.Font.Size = 10

# Column widths to fit the new content.
$ws2.Columns.Item(1).ColumnWidth = 12
$ws2.Columns.Item(2).ColumnWidth = 52

$ws2.Range("B14").Select()
